# Timesheet update 7.24
# Updates the employee ID on row 7 and the day-status letter codes
# (W/A/H/O/S/V/R) across rows 7-10, then leaves the selection on AA8
# to match the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - employee id
$ws.Range("B7").Value = 202207237495

# Row 7 day codes
$ws.Range("G7").Value  = "W"
$ws.Range("H7").Value  = "W"
$ws.Range("I7").Value  = "O"
$ws.Range("J7").Value  = "O"
$ws.Range("L7").Value  = "R"
$ws.Range("N7").Value  = "W"
$ws.Range("O7").Value  = "W"
$ws.Range("P7").Value  = "O"
$ws.Range("Q7").Value  = "O"
$ws.Range("R7").Value  = "A"
$ws.Range("U7").Value  = "S"
$ws.Range("V7").Value  = "W"
$ws.Range("W7").Value  = "O"
$ws.Range("X7").Value  = "O"
$ws.Range("AB7").Value = "V"
$ws.Range("AC7").Value = "V"
$ws.Range("AD7").Value = "O"
$ws.Range("AE7").Value = "O"
$ws.Range("AI7").Value = "W"

# Row 8 day codes
$ws.Range("G8").Value  = "W"
$ws.Range("H8").Value  = "W"
$ws.Range("I8").Value  = "O"
$ws.Range("J8").Value  = "O"
$ws.Range("N8").Value  = "W"
$ws.Range("O8").Value  = "W"
$ws.Range("P8").Value  = "O"
$ws.Range("Q8").Value  = "O"
$ws.Range("U8").Value  = "W"
$ws.Range("V8").Value  = "W"
$ws.Range("W8").Value  = "O"
$ws.Range("X8").Value  = "O"
$ws.Range("AB8").Value = "W"
$ws.Range("AC8").Value = "W"
$ws.Range("AD8").Value = "O"
$ws.Range("AE8").Value = "O"
$ws.Range("AI8").Value = "W"

# Row 9 day codes
$ws.Range("G9").Value  = "W"
$ws.Range("H9").Value  = "W"
$ws.Range("I9").Value  = "O"
$ws.Range("J9").Value  = "O"
$ws.Range("N9").Value  = "W"
$ws.Range("O9").Value  = "W"
$ws.Range("P9").Value  = "O"
$ws.Range("Q9").Value  = "O"
$ws.Range("U9").Value  = "W"
$ws.Range("V9").Value  = "W"
$ws.Range("W9").Value  = "O"
$ws.Range("X9").Value  = "O"
$ws.Range("AB9").Value = "W"
$ws.Range("AC9").Value = "W"
$ws.Range("AD9").Value = "O"
$ws.Range("AE9").Value = "O"
$ws.Range("AI9").Value = "W"

# Row 10 day codes
$ws.Range("G10").Value  = "W"
$ws.Range("H10").Value  = "W"
$ws.Range("I10").Value  = "O"
$ws.Range("J10").Value  = "O"
$ws.Range("N10").Value  = "W"
$ws.Range("O10").Value  = "W"
$ws.Range("P10").Value  = "O"
$ws.Range("Q10").Value  = "O"
$ws.Range("U10").Value  = "W"
$ws.Range("V10").Value  = "W"
$ws.Range("W10").Value  = "O"
$ws.Range("X10").Value  = "O"
$ws.Range("AB10").Value = "W"
$ws.Range("AC10").Value = "W"
$ws.Range("AD10").Value = "O"
$ws.Range("AE10").Value = "O"
$ws.Range("AI10").Value = "W"

# Final selection matches the saved workbook state
$ws.Range("AA8").Select()
